$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 45.5
$ws.Range("I5").Value = 47.77778
$ws.Range("J5").Value = 41.4
$ws.Range("K5").Value = 47.77778
$ws.Range("L5").Value = 41.4
$ws.Range("M5").Value = 67.22221999999999
$ws.Range("N5").Value = -271.4
# Row 51
$ws.Range("H51").Value = 10105332
$ws.Range("I51").Value = 22732772
$ws.Range("J51").Value = 3380
$ws.Range("K51").Value = 22732772
$ws.Range("L51").Value = 3380
$ws.Range("M51").Value = -22732288
$ws.Range("N51").Value = -4348
# Row 116
$ws.Range("H116").Value = 2824.2632
$ws.Range("I116").Value = 2577
$ws.Range("J116").Value = 3751.5
$ws.Range("K116").Value = 2577
$ws.Range("L116").Value = 3751.5
$ws.Range("M116").Value = 865
$ws.Range("N116").Value = -10635.5
# Row 135
$ws.Range("H135").Value = 1082.8158
$ws.Range("I135").Value = 782.40625
$ws.Range("J135").Value = 2685
$ws.Range("K135").Value = 7041.65625
$ws.Range("L135").Value = 24165
$ws.Range("M135").Value = -4506.65625
$ws.Range("N135").Value = -29235

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 28
$ws.Range("H28").Value = 12401.5625
$ws.Range("I28").Value = 6120.6665
$ws.Range("J28").Value = 20477
$ws.Range("K28").Value = 6120.6665
$ws.Range("L28").Value = 20477
$ws.Range("M28").Value = -5928.6665
# Row 41
$ws.Range("H41").Value = 2000
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 5000
$ws.Range("K41").Value = 1000
$ws.Range("L41").Value = 5000
$ws.Range("M41").Value = -586
$ws.Range("N41").Value = -5828
# Row 99
$ws.Range("H99").Value = 12401.5625
$ws.Range("I99").Value = 6120.6665
$ws.Range("J99").Value = 20477
$ws.Range("K99").Value = 6120.6665
$ws.Range("L99").Value = 20477
$ws.Range("M99").Value = -3125.6665
# Row 131
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
# Row 132
$ws.Range("H132").Value = 6613.383
$ws.Range("I132").Value = 6834.4585
$ws.Range("J132").Value = 6382.696
$ws.Range("K132").Value = 20503.3755
$ws.Range("L132").Value = 19148.088
$ws.Range("M132").Value = -17973.3755
$ws.Range("N132").Value = -24208.088
# Row 139
$ws.Range("H139").Value = 34575
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 34575
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 34575
$ws.Range("N139").Value = -44855

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 25
$ws.Range("H25").Value = 757
$ws.Range("I25").Value = 757
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 757
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -522
# Row 35
$ws.Range("H35").Value = 70000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 70000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 70000
$ws.Range("N35").Value = -70620

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 490
$ws.Range("I17").Value = 490
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 490
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -316
# Row 31
$ws.Range("H31").Value = 2273.842
$ws.Range("I31").Value = 1391.6716
$ws.Range("J31").Value = 4384.75
$ws.Range("K31").Value = 1391.6716
$ws.Range("L31").Value = 4384.75
$ws.Range("M31").Value = -1096.6716
$ws.Range("N31").Value = -4974.75
# Row 34
$ws.Range("H34").Value = 2273.842
$ws.Range("I34").Value = 1391.6716
$ws.Range("J34").Value = 4384.75
$ws.Range("K34").Value = 1391.6716
$ws.Range("L34").Value = 4384.75
$ws.Range("M34").Value = -1189.6716
$ws.Range("N34").Value = -4788.75
# Row 36
$ws.Range("H36").Value = 4024
$ws.Range("I36").Value = 4024
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 4024
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -3636
# Row 40
$ws.Range("H40").Value = 4024
$ws.Range("I40").Value = 4024
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4024
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3864
# Row 58
$ws.Range("H58").Value = 5085.4614
$ws.Range("I58").Value = 2375.7222
$ws.Range("J58").Value = 11182.375
$ws.Range("K58").Value = 2375.7222
$ws.Range("L58").Value = 11182.375
$ws.Range("M58").Value = -2172.7222
$ws.Range("N58").Value = -11588.375
# Row 88
$ws.Range("H88").Value = 23298.889
$ws.Range("I88").Value = 15000
$ws.Range("J88").Value = 24336.25
$ws.Range("K88").Value = 15000
$ws.Range("L88").Value = 24336.25
$ws.Range("M88").Value = -14594
$ws.Range("N88").Value = -25148.25
# Row 91
$ws.Range("H91").Value = 23298.889
$ws.Range("I91").Value = 15000
$ws.Range("J91").Value = 24336.25
$ws.Range("K91").Value = 15000
$ws.Range("L91").Value = 24336.25
$ws.Range("M91").Value = -13596
$ws.Range("N91").Value = -27144.25
# Row 136
$ws.Range("H136").Value = 5085.4614
$ws.Range("I136").Value = 2375.7222
$ws.Range("J136").Value = 11182.375
$ws.Range("K136").Value = 7127.1666
$ws.Range("L136").Value = 33547.125
$ws.Range("M136").Value = -4577.1666
$ws.Range("N136").Value = -38647.125

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 1809.3231
$ws.Range("I122").Value = 465.86365
$ws.Range("J122").Value = 2496.6743
$ws.Range("K122").Value = 4192.77285
$ws.Range("L122").Value = 22470.0687
$ws.Range("M122").Value = -1742.77285
$ws.Range("N122").Value = -27370.0687

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 2159.4
$ws.Range("I9").Value = 929.6667
$ws.Range("J9").Value = 4004
$ws.Range("K9").Value = 929.6667
$ws.Range("L9").Value = 4004
$ws.Range("M9").Value = -759.6667
# Row 31
$ws.Range("H31").Value = 1232.75
$ws.Range("I31").Value = 1232.75
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1232.75
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -940.75
# Row 37
$ws.Range("H37").Value = 1232.75
$ws.Range("I37").Value = 1232.75
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1232.75
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -955.75
# Row 80
$ws.Range("H80").Value = 2412.24
$ws.Range("I80").Value = 1957.1428
$ws.Range("J80").Value = 2589.2222
$ws.Range("K80").Value = 1957.1428
$ws.Range("L80").Value = 2589.2222
$ws.Range("M80").Value = -959.1428000000001
$ws.Range("N80").Value = -4585.2222
# Row 83
$ws.Range("H83").Value = 2412.24
$ws.Range("I83").Value = 1957.1428
$ws.Range("J83").Value = 2589.2222
$ws.Range("K83").Value = 9785.714
$ws.Range("L83").Value = 12946.111
$ws.Range("M83").Value = -4793.714
$ws.Range("N83").Value = -22930.111

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 558
$ws.Range("I9").Value = 447.5
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 447.5
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = -223.5
# Row 30
$ws.Range("H30").Value = 30329.334
$ws.Range("I30").Value = 640.6667
$ws.Range("J30").Value = 60018
$ws.Range("K30").Value = 640.6667
$ws.Range("L30").Value = 60018
$ws.Range("M30").Value = -532.6667
# Row 35
$ws.Range("H35").Value = 517.875
$ws.Range("I35").Value = 517.875
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 517.875
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -181.875
# Row 46
$ws.Range("H46").Value = 385206.47
$ws.Range("I46").Value = 486.33334
$ws.Range("J46").Value = 909824.8
$ws.Range("K46").Value = 486.33334
$ws.Range("L46").Value = 909824.8
$ws.Range("M46").Value = -298.33334
# Row 93
$ws.Range("H93").Value = 1779
$ws.Range("I93").Value = 1415.1428
$ws.Range("J93").Value = 2142.8572
$ws.Range("K93").Value = 1415.1428
$ws.Range("L93").Value = 2142.8572
$ws.Range("M93").Value = -167.1428000000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 44
$ws.Range("H44").Value = 20000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 20000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -21108
# Row 52
$ws.Range("H52").Value = 5750
$ws.Range("I52").Value = 5500
$ws.Range("J52").Value = 6000
$ws.Range("K52").Value = 5500
$ws.Range("L52").Value = 6000
$ws.Range("M52").Value = -5274
$ws.Range("N52").Value = -6452
# Row 69
$ws.Range("H69").Value = 12411.429
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 12411.429
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 12411.429
$ws.Range("N69").Value = -13909.429
# Row 72
$ws.Range("H72").Value = 12411.429
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 12411.429
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 37234.287
$ws.Range("N72").Value = -44722.287
# Row 100
$ws.Range("H100").Value = 845.3889
$ws.Range("I100").Value = 466.16666
$ws.Range("J100").Value = 1603.8334
$ws.Range("K100").Value = 932.33332
$ws.Range("L100").Value = 3207.6668
$ws.Range("M100").Value = -391.33332
$ws.Range("N100").Value = -4289.6668
